$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 22.3
$ws.Range("B3").Value = "28/01/2025"
$ws.Range("C3").Value = "Mens Casual Premium Slim Fit T-Shirts "
